# Reformat the "mean ± stdev" metric cells (and the one plain p-value cell)
# in the "Peptide Evaluation" results table down to 2 decimal places, e.g.
# "0.8051 ± 0.0871" -> "0.81 ± 0.09". Any trailing significance marker
# (" *") is preserved as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "0.81 ± 0.09"
$ws.Range("D2").Value = "0.68 ± 0.17"
$ws.Range("E2").Value = "0.69 ± 0.11"
$ws.Range("F2").Value = "0.56 ± 0.26"
$ws.Range("G2").Value = "0.60 ± 0.21"
$ws.Range("H2").Value = "0.47 ± 0.22"
$ws.Range("C3").Value = "0.72 ± 0.03"
$ws.Range("D3").Value = "0.53 ± 0.06"
$ws.Range("E3").Value = "0.63 ± 0.03"
$ws.Range("F3").Value = "0.31 ± 0.07"
$ws.Range("G3").Value = "0.41 ± 0.07"
$ws.Range("H3").Value = "0.28 ± 0.05"
$ws.Range("C4").Value = "0.52 ± 0.03"
$ws.Range("D4").Value = "0.00 ± 0.00"
$ws.Range("E4").Value = "0.00 ± 0.00"
$ws.Range("F4").Value = "0.00 ± 0.00"
$ws.Range("G4").Value = "0.00 ± 0.00"
$ws.Range("H4").Value = "0.00 ± 0.00"
$ws.Range("C5").Value = "0.57 ± 0.04"
$ws.Range("D5").Value = "0.01 ± 0.02"
$ws.Range("E5").Value = "0.10 ± 0.30"
$ws.Range("F5").Value = "0.00 ± 0.00"
$ws.Range("G5").Value = "0.00 ± 0.00"
$ws.Range("H5").Value = "0.01 ± 0.02"
$ws.Range("C6").Value = "0.73 ± 0.02"
$ws.Range("D6").Value = "0.57 ± 0.04"
$ws.Range("E6").Value = "0.64 ± 0.03"
$ws.Range("F6").Value = "0.36 ± 0.05"
$ws.Range("G6").Value = "0.46 ± 0.04"
$ws.Range("H6").Value = "0.31 ± 0.04"
$ws.Range("C8").Value = "0.86 ± 0.08"
$ws.Range("D8").Value = "0.77 ± 0.13"
$ws.Range("E8").Value = "0.70 ± 0.13"
$ws.Range("F8").Value = "0.73 ± 0.20"
$ws.Range("G8").Value = "0.70 ± 0.16"
$ws.Range("H8").Value = "0.57 ± 0.21"
$ws.Range("C9").Value = "0.93 ± 0.02"
$ws.Range("D9").Value = "0.87 ± 0.02"
$ws.Range("E9").Value = "0.88 ± 0.02"
$ws.Range("F9").Value = "0.81 ± 0.04"
$ws.Range("G9").Value = "0.84 ± 0.03"
$ws.Range("H9").Value = "0.77 ± 0.04"
$ws.Range("C10").Value = "0.59 ± 0.11"
$ws.Range("D10").Value = "0.09 ± 0.27"
$ws.Range("E10").Value = "0.09 ± 0.26"
$ws.Range("F10").Value = "0.08 ± 0.25"
$ws.Range("G10").Value = "0.08 ± 0.25"
$ws.Range("H10").Value = "0.08 ± 0.23"
$ws.Range("C11").Value = "0.62 ± 0.06"
$ws.Range("D11").Value = "0.10 ± 0.07"
$ws.Range("E11").Value = "0.36 ± 0.30"
$ws.Range("F11").Value = "0.01 ± 0.01"
$ws.Range("G11").Value = "0.03 ± 0.02"
$ws.Range("H11").Value = "0.02 ± 0.05"
$ws.Range("C12").Value = "0.93 ± 0.02"
$ws.Range("D12").Value = "0.88 ± 0.02"
$ws.Range("E12").Value = "0.88 ± 0.02"
$ws.Range("F12").Value = "0.82 ± 0.05"
$ws.Range("G12").Value = "0.85 ± 0.03"
$ws.Range("H12").Value = "0.78 ± 0.04"
$ws.Range("C14").Value = "0.82 ± 0.12"
$ws.Range("D14").Value = "0.70 ± 0.25"
$ws.Range("E14").Value = "0.68 ± 0.17"
$ws.Range("F14").Value = "0.68 ± 0.32"
$ws.Range("G14").Value = "0.64 ± 0.29"
$ws.Range("H14").Value = "0.53 ± 0.28"
$ws.Range("C15").Value = "0.94 ± 0.02"
$ws.Range("D15").Value = "0.88 ± 0.02"
$ws.Range("E15").Value = "0.89 ± 0.03 *"
$ws.Range("F15").Value = "0.82 ± 0.04"
$ws.Range("G15").Value = "0.85 ± 0.02"
$ws.Range("H15").Value = "0.78 ± 0.03 *"
$ws.Range("C16").Value = "0.79 ± 0.02"
$ws.Range("D16").Value = "0.63 ± 0.04"
$ws.Range("E16").Value = "0.85 ± 0.04"
$ws.Range("F16").Value = "0.41 ± 0.06"
$ws.Range("G16").Value = "0.55 ± 0.05"
$ws.Range("H16").Value = "0.48 ± 0.03"
$ws.Range("C17").Value = "0.87 ± 0.02"
$ws.Range("D17").Value = "0.76 ± 0.03"
$ws.Range("E17").Value = "0.86 ± 0.02"
$ws.Range("F17").Value = "0.61 ± 0.06"
$ws.Range("G17").Value = "0.71 ± 0.04"
$ws.Range("H17").Value = "0.62 ± 0.04"
$ws.Range("C18").Value = "0.94 ± 0.02 *"
$ws.Range("D18").Value = "0.88 ± 0.02"
$ws.Range("E18").Value = "0.87 ± 0.02"
$ws.Range("F18").Value = "0.82 ± 0.03"
$ws.Range("G18").Value = "0.85 ± 0.02"
$ws.Range("H18").Value = "0.78 ± 0.03"
# "0.00" reads as a plain number, so a normal Value assignment would be
# coerced to the number 0 instead of staying text. Type it the way a user
# would to force text (leading apostrophe), then clear the resulting
# quote-prefix formatting so the cell's appearance/style is unaffected.
$ws.Range("E19").Value = "'0.00"
$ws.Range("E19").ClearFormats()
